# Applies the data-refresh edit described by the commit "Add files via upload":
#  - updates a batch of measured values on "Extended Fig3a" (columns B and D)
#  - updates one measured value on "Extended Fig3b" (E10)
#  - moves the active-cell selection on both sheets to match where the author
#    was last working
#  - shrinks the saved window size to match the author's window

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Extended Fig3a")
$ws2 = $wb.Worksheets.Item("Extended Fig3b")

# --- Extended Fig3a: refreshed measurements in columns B (row) and D (row) ---
$ws1.Cells.Item(4, 2).Value = 13.2
$ws1.Cells.Item(4, 4).Value = 12
$ws1.Cells.Item(5, 2).Value = 13.04
$ws1.Cells.Item(5, 4).Value = 12.84
$ws1.Cells.Item(6, 2).Value = 13.08
$ws1.Cells.Item(6, 4).Value = 14.36
$ws1.Cells.Item(7, 2).Value = 12.32
$ws1.Cells.Item(7, 4).Value = 13.82
$ws1.Cells.Item(8, 2).Value = 11.3
$ws1.Cells.Item(8, 4).Value = 13.52
$ws1.Cells.Item(9, 2).Value = 11.89
$ws1.Cells.Item(9, 4).Value = 13.78
$ws1.Cells.Item(10, 2).Value = 13.72
$ws1.Cells.Item(10, 4).Value = 15.71
$ws1.Cells.Item(11, 2).Value = 14.69
$ws1.Cells.Item(11, 4).Value = 16.16
$ws1.Cells.Item(12, 2).Value = 13.65
$ws1.Cells.Item(12, 4).Value = 13.36
$ws1.Cells.Item(13, 2).Value = 12.43
$ws1.Cells.Item(13, 4).Value = 9.83
$ws1.Cells.Item(14, 2).Value = 10.68
$ws1.Cells.Item(14, 4).Value = 8.32
$ws1.Cells.Item(15, 2).Value = 8.13
$ws1.Cells.Item(15, 4).Value = 6.01
$ws1.Cells.Item(16, 2).Value = 7.72
$ws1.Cells.Item(16, 4).Value = 3.7
$ws1.Cells.Item(17, 2).Value = 7.11
$ws1.Cells.Item(17, 4).Value = 4.04
$ws1.Cells.Item(18, 2).Value = 6.61
$ws1.Cells.Item(18, 4).Value = 4.2
$ws1.Cells.Item(19, 2).Value = 7
$ws1.Cells.Item(19, 4).Value = 4.8
$ws1.Cells.Item(20, 2).Value = 7.52
$ws1.Cells.Item(20, 4).Value = 5.42
$ws1.Cells.Item(21, 2).Value = 8.67
$ws1.Cells.Item(21, 4).Value = 7.06
$ws1.Cells.Item(22, 2).Value = 9.62
$ws1.Cells.Item(22, 4).Value = 9.35
$ws1.Cells.Item(23, 2).Value = 11.29
$ws1.Cells.Item(23, 4).Value = 11.02
$ws1.Cells.Item(24, 2).Value = 12.64
$ws1.Cells.Item(24, 4).Value = 10.98
$ws1.Cells.Item(25, 2).Value = 13.65
$ws1.Cells.Item(25, 4).Value = 10.52
$ws1.Cells.Item(26, 2).Value = 13.66
$ws1.Cells.Item(26, 4).Value = 11.29
$ws1.Cells.Item(27, 2).Value = 13.25
$ws1.Cells.Item(27, 4).Value = 10.94

# --- Extended Fig3b: one refreshed measurement ---
$ws2.Cells.Item(10, 5).Value = 4.17

# --- Restore the author's on-screen selection on each sheet ---
$ws1.Activate()
$ws1.Range("D21").Select()

$ws2.Activate()
$ws2.Range("F31").Select()

# Leave the workbook focused back on the first sheet, matching tabSelected="1"
$ws1.Activate()

# --- Match the saved window size in the workbook view ---
$win = $excel.ActiveWindow
$win.Width = 15960
$win.Height = 9255
